$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 124, shifting existing rows 124:144 down to 125:145
$ws.Rows(124).Insert()

# Populate the newly inserted row 124 with the new data record
$ws.Cells.Item(124, 1).Value = 4
$ws.Cells.Item(124, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(124, 3).Value = "Los Lagos"
$ws.Cells.Item(124, 4).Value = 44504
$ws.Cells.Item(124, 5).Value = 10
$ws.Cells.Item(124, 6).Value = 100112032
$ws.Cells.Item(124, 7).Value = "Zapallo italiano"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 80
$ws.Cells.Item(124, 11).Value = 11000
$ws.Cells.Item(124, 12).Value = 11000
$ws.Cells.Item(124, 13).Value = 11000
$ws.Cells.Item(124, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(124, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(124, 16).Value = 220
$ws.Cells.Item(124, 17).Value = 50
$ws.Cells.Item(124, 18).Value = "Hortaliza"
